$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at J:K (pushes old Modality..filename columns
# from J..P to L..R), then insert one more blank column at M (pushes them to
# their final resting place L, N, O, P, Q, R, S).
$ws.Columns("J:K").Insert()
$ws.Columns("M:M").Insert()

# Populate the three new header cells introduced by this change.
$ws.Range("J1").Value = "Manufacturer"
$ws.Range("K1").Value = "ManufacturerModelName"
$ws.Range("M1").Value = "StationName"

# Update the active selection to match the edited workbook.
$null = $ws.Range("M1").Select()
